# eventsliders.xlsx edit:
#  - Commands sheet: insert a new row 58 for the new YOCTOPUCE "yset" relay
#    IO command (between "off(c[,sn])" and "flip(c[,sn])"), shifting all
#    following rows down by one.
#  - Update the active selection on the Commands sheet to track the moved
#    row (was B65:C65, now B58:C58 after the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at 58, pushing "flip(c[,sn])" (and everything after it)
# down by one row.
$ws.Rows("58:58").Insert()

# Populate the new row with the yset(c,b[,sn]) command documentation.
$ws.Range("B58").Value = "yset(c,b[,sn])"
$ws.Range("C58").Value = "YOCTOPUCE Relay Output: switches channel c of the relay module off (b=0) and on (b=1)"

# Match the row height used by the other command rows in this table.
$ws.Rows("58:58").RowHeight = 13.8

# Refresh the active selection to the newly inserted row.
$ws.Range("B58:C58").Select() | Out-Null
